# Apply scheduled market-price data updates to each Leve profit sheet.
# Values below come from the latest Universalis price snapshot pull.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 34040.8
$ws.Range("J3").Value = 34040.8
$ws.Range("L3").Value = 34040.8
$ws.Range("N3").Value = -34268.8
$ws.Range("H17").Value = 2570351
$ws.Range("J17").Value = 2763090
$ws.Range("L17").Value = 8289270
$ws.Range("N17").Value = -8289606
$ws.Range("H101").Value = 2082
$ws.Range("I101").Value = 2193.7144
$ws.Range("J101").Value = 1300
$ws.Range("K101").Value = 6581.1432
$ws.Range("L101").Value = 3900
$ws.Range("M101").Value = -4959.1432
$ws.Range("N101").Value = -7144
$ws.Range("H102").Value = 34040.8
$ws.Range("J102").Value = 34040.8
$ws.Range("L102").Value = 34040.8
$ws.Range("N102").Value = -40530.8
$ws.Range("H138").Value = 4833299
$ws.Range("I138").Value = 8548499
$ws.Range("J138").Value = 3537.8667
$ws.Range("K138").Value = 25645497
$ws.Range("L138").Value = 10613.6001
$ws.Range("M138").Value = -25640357
$ws.Range("N138").Value = -20893.6001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2707.625
$ws.Range("I61").Value = 2415.2163
$ws.Range("J61").Value = 6314
$ws.Range("K61").Value = 2415.2163
$ws.Range("L61").Value = 6314
$ws.Range("M61").Value = -2203.2163
$ws.Range("N61").Value = -6738
$ws.Range("H74").Value = 1459.1212
$ws.Range("I74").Value = 1133.0416
$ws.Range("K74").Value = 1133.0416
$ws.Range("M74").Value = -259.0416
$ws.Range("H77").Value = 1459.1212
$ws.Range("I77").Value = 1133.0416
$ws.Range("K77").Value = 5665.208000000001
$ws.Range("M77").Value = -1297.208000000001
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 35184.668
$ws.Range("J109").Value = 35184.668
$ws.Range("L109").Value = 35184.668
$ws.Range("N109").Value = -37958.668
$ws.Range("H136").Value = 2707.625
$ws.Range("I136").Value = 2415.2163
$ws.Range("J136").Value = 6314
$ws.Range("K136").Value = 7245.6489
$ws.Range("L136").Value = 18942
$ws.Range("M136").Value = -4695.6489
$ws.Range("N136").Value = -24042

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 28675.268
$ws.Range("J132").Value = 28675.268
$ws.Range("L132").Value = 28675.268
$ws.Range("N132").Value = -38795.268
$ws.Range("H134").Value = 3213.1155
$ws.Range("I134").Value = 1891.7097
$ws.Range("J134").Value = 5163.7617
$ws.Range("K134").Value = 5675.1291
$ws.Range("L134").Value = 15491.2851
$ws.Range("M134").Value = -3140.1291
$ws.Range("N134").Value = -20561.2851

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 94628.57000000001
$ws.Range("J37").Value = 94628.57000000001
$ws.Range("L37").Value = 283885.71
$ws.Range("N37").Value = -284109.71
$ws.Range("H70").Value = 3721.875
$ws.Range("I70").Value = 1487.5
$ws.Range("K70").Value = 4462.5
$ws.Range("M70").Value = -4147.5
$ws.Range("H73").Value = 3721.875
$ws.Range("I73").Value = 1487.5
$ws.Range("K73").Value = 4462.5
$ws.Range("M73").Value = -3370.5
$ws.Range("H98").Value = 1183.72
$ws.Range("I98").Value = 156.14285
$ws.Range("J98").Value = 1583.3334
$ws.Range("K98").Value = 468.42855
$ws.Range("L98").Value = 4750.0002
$ws.Range("M98").Value = 1029.57145
$ws.Range("N98").Value = -7746.0002

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 15000
$ws.Range("J64").Value = 15000
$ws.Range("L64").Value = 15000
$ws.Range("H67").Value = 15000
$ws.Range("J67").Value = 15000
$ws.Range("L67").Value = 15000
$ws.Range("H132").Value = 6355.174
$ws.Range("I132").Value = 7843
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 23529
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -20999
$ws.Range("N132").Value = -8057

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1230.6471
$ws.Range("I16").Value = 1416.6154
$ws.Range("J16").Value = 626.25
$ws.Range("K16").Value = 1416.6154
$ws.Range("L16").Value = 626.25
$ws.Range("M16").Value = -1246.6154
$ws.Range("N16").Value = -966.25
$ws.Range("H93").Value = 932.1070999999999
$ws.Range("I93").Value = 871.96
$ws.Range("J93").Value = 1433.3334
$ws.Range("K93").Value = 871.96
$ws.Range("L93").Value = 1433.3334
$ws.Range("M93").Value = 376.04
$ws.Range("N93").Value = -3929.3334
$ws.Range("H100").Value = 2290
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
$ws.Range("H102").Value = 50256.4
$ws.Range("J102").Value = 50256.4
$ws.Range("L102").Value = 50256.4
$ws.Range("N102").Value = -56746.4
$ws.Range("H136").Value = 4279.3335
$ws.Range("I136").Value = 1993.1538
$ws.Range("J136").Value = 33999.668
$ws.Range("K136").Value = 5979.4614
$ws.Range("L136").Value = 101999.004
$ws.Range("M136").Value = -3429.4614
$ws.Range("N136").Value = -107099.004

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3456.125
$ws.Range("I96").Value = 2225
$ws.Range("J96").Value = 3866.5
$ws.Range("K96").Value = 2225
$ws.Range("L96").Value = 3866.5
$ws.Range("M96").Value = -852
$ws.Range("N96").Value = -6612.5
$ws.Range("H100").Value = 1900.6
$ws.Range("I100").Value = 2450
$ws.Range("J100").Value = 1534.3334
$ws.Range("K100").Value = 4900
$ws.Range("L100").Value = 3068.6668
$ws.Range("M100").Value = -4359
$ws.Range("N100").Value = -4150.6668
$ws.Range("H106").Value = 43000
$ws.Range("J106").Value = 43000
$ws.Range("L106").Value = 43000
$ws.Range("H109").Value = 33377
$ws.Range("J109").Value = 33377
$ws.Range("L109").Value = 33377
$ws.Range("N109").Value = -36151
$ws.Range("H136").Value = 1615.9736
$ws.Range("I136").Value = 951.9655
$ws.Range("J136").Value = 3755.5557
$ws.Range("K136").Value = 2855.8965
$ws.Range("L136").Value = 11266.6671
$ws.Range("M136").Value = -305.8964999999998
$ws.Range("N136").Value = -16366.6671
